$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two more "Mặt sau" rows (9 and 10), identical in content to the
# existing rows 4-8, extending the used range from A1:I8 to A1:I10.
$values = @("Mặt sau", "PHAM DUY LONG", "S Trà Co, Thanh Cái, Qung NInh phó Móng Khu Trang Ginl Trà Co, Thanh Móng Cál, phó", "03/12/2006", "022206004066", "0v12/2031", "Việt Nam", "Hải Xuan, Thành phố Móng Cái, Quảng Ninh Hải Xuán, Thành phó Móng Cá", "Nam")

foreach ($r in 9, 10) {
    for ($c = 1; $c -le 9; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Style = "Normal"
        $cell.NumberFormat = "@"
        $cell.Value = $values[$c - 1]
        $cell.Style = "Normal"
    }
}
